$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers (values first, then copy the existing header's formatting over)
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)

# New data columns
$data = @(
    @(91.51563705136908, 237743, 317.8382352941176),
    @(89.34822657184763, 25810, 314.7560975609756),
    @(89.68836503732236, 168645, 143.0407124681934),
    @(92.47271837489988, 35001, 171.5735294117647),
    @(18.1379821112245, 2042, 13.89115646258503),
    @(37.04336315460473, 166, 33.2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $data[$i][0]
    $ws.Cells.Item($row, 13).Value = $data[$i][1]
    $ws.Cells.Item($row, 14).Value = $data[$i][2]
}
